# Scheduled market-price refresh: update currentAveragePrice* / Leve*Price* /
# LeveProfit* columns (H:N) for the affected rows across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 14900.869
$ws.Range("I40").Value = 26689.6
$ws.Range("J40").Value = 11626.223
$ws.Range("K40").Value = 26689.6
$ws.Range("L40").Value = 11626.223
$ws.Range("M40").Value = -26514.6
$ws.Range("N40").Value = -11976.223
$ws.Range("H55").Value = 621.625
$ws.Range("I55").Value = 754.6
$ws.Range("K55").Value = 754.6
$ws.Range("M55").Value = -540.6
$ws.Range("H127").Value = 2743.8857
$ws.Range("I127").Value = 1157.3846
$ws.Range("J127").Value = 3681.3635
$ws.Range("K127").Value = 3472.1538
$ws.Range("L127").Value = 11044.0905
$ws.Range("M127").Value = 1487.8462
$ws.Range("N127").Value = -20964.0905
$ws.Range("H132").Value = 8600.431
$ws.Range("I132").Value = 1935.1515
$ws.Range("J132").Value = 13382.044
$ws.Range("K132").Value = 5805.4545
$ws.Range("L132").Value = 40146.132
$ws.Range("M132").Value = -3275.4545
$ws.Range("N132").Value = -45206.132
$ws.Range("H133").Value = 93696.664
$ws.Range("J133").Value = 93696.664
$ws.Range("L133").Value = 93696.664
$ws.Range("N133").Value = -103816.664
$ws.Range("H137").Value = 4542.476
$ws.Range("J137").Value = 6366.5
$ws.Range("L137").Value = 19099.5
$ws.Range("N137").Value = -24199.5
$ws.Range("H138").Value = 6557.4443
$ws.Range("J138").Value = 7463.773
$ws.Range("L138").Value = 22391.319
$ws.Range("N138").Value = -32671.319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4328.467
$ws.Range("I32").Value = 2456.2917
$ws.Range("J32").Value = 11817.167
$ws.Range("K32").Value = 2456.2917
$ws.Range("L32").Value = 11817.167
$ws.Range("M32").Value = -2169.2917
$ws.Range("N32").Value = -12391.167
$ws.Range("H110").Value = 821000.75
$ws.Range("I110").Value = 1078542.5
$ws.Range("J110").Value = 5451.8335
$ws.Range("K110").Value = 1078542.5
$ws.Range("L110").Value = 5451.8335
$ws.Range("M110").Value = -1076497.5
$ws.Range("N110").Value = -9541.833500000001
$ws.Range("H122").Value = 15156421
$ws.Range("I122").Value = 23810662
$ws.Range("J122").Value = 11499.75
$ws.Range("K122").Value = 71431986
$ws.Range("L122").Value = 34499.25
$ws.Range("M122").Value = -71429536
$ws.Range("N122").Value = -39399.25
$ws.Range("H132").Value = 31139.727
$ws.Range("I132").Value = 36881.875
$ws.Range("K132").Value = 110645.625
$ws.Range("M132").Value = -108115.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 386.5
$ws.Range("I22").Value = 458.125
$ws.Range("K22").Value = 458.125
$ws.Range("M22").Value = -285.125
$ws.Range("H86").Value = 27779648
$ws.Range("I86").Value = 1727.3182
$ws.Range("J86").Value = 71430664
$ws.Range("K86").Value = 1727.3182
$ws.Range("L86").Value = 71430664
$ws.Range("M86").Value = -604.3181999999999
$ws.Range("N86").Value = -71432910
$ws.Range("H89").Value = 27779648
$ws.Range("I89").Value = 1727.3182
$ws.Range("J89").Value = 71430664
$ws.Range("K89").Value = 8636.591
$ws.Range("L89").Value = 357153320
$ws.Range("M89").Value = -3020.591
$ws.Range("N89").Value = -357164552
$ws.Range("H108").Value = 91999.5
$ws.Range("J108").Value = 91999.5
$ws.Range("L108").Value = 91999.5
$ws.Range("N108").Value = -99679.5
$ws.Range("H134").Value = 3254.6667
$ws.Range("I134").Value = 3750
$ws.Range("K134").Value = 11250
$ws.Range("M134").Value = -8715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1429.6666
$ws.Range("I16").Value = 1345.8334
$ws.Range("K16").Value = 1345.8334
$ws.Range("M16").Value = -1058.8334
$ws.Range("H22").Value = 178.6923
$ws.Range("I22").Value = 173.6
$ws.Range("J22").Value = 181.875
$ws.Range("K22").Value = 173.6
$ws.Range("L22").Value = 181.875
$ws.Range("M22").Value = 176.4
$ws.Range("N22").Value = -881.875
$ws.Range("H31").Value = 3309.0728
$ws.Range("J31").Value = 7312.375
$ws.Range("L31").Value = 7312.375
$ws.Range("N31").Value = -7902.375
$ws.Range("H34").Value = 3309.0728
$ws.Range("J34").Value = 7312.375
$ws.Range("L34").Value = 7312.375
$ws.Range("N34").Value = -7716.375
$ws.Range("H107").Value = 729334.8
$ws.Range("I107").Value = 1212396.8
$ws.Range("J107").Value = 4741.8
$ws.Range("K107").Value = 1212396.8
$ws.Range("L107").Value = 4741.8
$ws.Range("M107").Value = -1210476.8
$ws.Range("N107").Value = -8581.799999999999
$ws.Range("H113").Value = 1429.6666
$ws.Range("I113").Value = 1345.8334
$ws.Range("K113").Value = 1345.8334
$ws.Range("M113").Value = 824.1666
$ws.Range("H122").Value = 3208510
$ws.Range("I122").Value = 4811068.5
$ws.Range("K122").Value = 14433205.5
$ws.Range("M122").Value = -14430755.5
$ws.Range("H132").Value = 2387.1428
$ws.Range("I132").Value = 2542
$ws.Range("K132").Value = 7626
$ws.Range("M132").Value = -5096
$ws.Range("H134").Value = 2798.5652
$ws.Range("I134").Value = 2335.3157
$ws.Range("K134").Value = 7005.9471
$ws.Range("M134").Value = -4470.9471

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36507576
$ws.Range("I4").Value = 1083634.4
$ws.Range("J4").Value = 178203330
$ws.Range("K4").Value = 3250903.2
$ws.Range("L4").Value = 534609990
$ws.Range("M4").Value = -3250791.2
$ws.Range("N4").Value = -534610214
$ws.Range("H6").Value = 202.5
$ws.Range("I6").Value = 202.5
$ws.Range("K6").Value = 607.5
$ws.Range("M6").Value = -494.5
$ws.Range("H12").Value = 224.34616
$ws.Range("I12").Value = 313.42856
$ws.Range("J12").Value = 191.52632
$ws.Range("K12").Value = 940.28568
$ws.Range("L12").Value = 574.5789600000001
$ws.Range("M12").Value = -767.28568
$ws.Range("N12").Value = -920.5789600000001
$ws.Range("H26").Value = 236.75
$ws.Range("I26").Value = 236.75
$ws.Range("K26").Value = 710.25
$ws.Range("M26").Value = -422.25
$ws.Range("H131").Value = 2108.0469
$ws.Range("I131").Value = 1714.5333
$ws.Range("J131").Value = 2228.5103
$ws.Range("K131").Value = 5143.5999
$ws.Range("L131").Value = 6685.5309
$ws.Range("M131").Value = -103.5999000000002
$ws.Range("N131").Value = -16765.5309

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1913209
$ws.Range("I70").Value = 3182348.8
$ws.Range("K70").Value = 3182348.8
$ws.Range("M70").Value = -3182078.8
$ws.Range("H73").Value = 1913209
$ws.Range("I73").Value = 3182348.8
$ws.Range("K73").Value = 3182348.8
$ws.Range("M73").Value = -3181412.8
$ws.Range("H102").Value = 5124.553
$ws.Range("J102").Value = 8074.25
$ws.Range("L102").Value = 8074.25
$ws.Range("N102").Value = -11318.25
$ws.Range("H122").Value = 3764533
$ws.Range("I122").Value = 5853938
$ws.Range("J122").Value = 3604.2
$ws.Range("K122").Value = 17561814
$ws.Range("L122").Value = 10812.6
$ws.Range("M122").Value = -17559364
$ws.Range("N122").Value = -15712.6
$ws.Range("H126").Value = 3932.9048
$ws.Range("I126").Value = 1976.0588
$ws.Range("J126").Value = 12249.5
$ws.Range("K126").Value = 5928.1764
$ws.Range("L126").Value = 36748.5
$ws.Range("M126").Value = -3458.1764
$ws.Range("N126").Value = -41688.5
$ws.Range("H132").Value = 4870.4253
$ws.Range("I132").Value = 4665.161
$ws.Range("K132").Value = 13995.483
$ws.Range("M132").Value = -11465.483

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4466.9165
$ws.Range("I7").Value = 3632.2368
$ws.Range("K7").Value = 3632.2368
$ws.Range("M7").Value = -3520.2368
$ws.Range("H40").Value = 7440.1665
$ws.Range("J40").Value = 7440.1665
$ws.Range("L40").Value = 7440.1665
$ws.Range("N40").Value = -7712.1665
$ws.Range("H46").Value = 5664.12
$ws.Range("J46").Value = 5994.6523
$ws.Range("L46").Value = 5994.6523
$ws.Range("N46").Value = -6370.6523
$ws.Range("H55").Value = 749.9167
$ws.Range("I55").Value = 727.25
$ws.Range("J55").Value = 795.25
$ws.Range("K55").Value = 727.25
$ws.Range("L55").Value = 795.25
$ws.Range("M55").Value = -554.25
$ws.Range("N55").Value = -1141.25
$ws.Range("H126").Value = 4466.9165
$ws.Range("I126").Value = 3632.2368
$ws.Range("K126").Value = 10896.7104
$ws.Range("M126").Value = -8426.7104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2774.2642
$ws.Range("I122").Value = 2670.9534
$ws.Range("J122").Value = 3218.5
$ws.Range("K122").Value = 8012.860199999999
$ws.Range("L122").Value = 9655.5
$ws.Range("M122").Value = -5562.860199999999
$ws.Range("N122").Value = -14555.5
$ws.Range("H126").Value = 1276.125
$ws.Range("I126").Value = 1297.125
$ws.Range("J126").Value = 1255.125
$ws.Range("K126").Value = 3891.375
$ws.Range("L126").Value = 3765.375
$ws.Range("M126").Value = -1421.375
$ws.Range("N126").Value = -8705.375
$ws.Range("H136").Value = 8548.037
$ws.Range("I136").Value = 3442.8667
$ws.Range("K136").Value = 10328.6001
$ws.Range("M136").Value = -7778.6001
